# This script applies the Gilgamesh_Profits update: recalculated leve
# market-price figures (currentAveragePrice / NQ / HQ, Leve prices, and
# profit columns) for a batch of rows across the ALC, ARM, BSM, CRP, CUL,
# GSM, LTW, and WVR sheets, per the scheduled pricing-data refresh.
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 923.4706
$ws.Range("J53").Value = 208.25
$ws.Range("L53").Value = 208.25
$ws.Range("N53").Value = -1482.25
$ws.Range("H55").Value = 538.2941
$ws.Range("J55").Value = 581.2857
$ws.Range("L55").Value = 581.2857
$ws.Range("N55").Value = -1009.2857
$ws.Range("H74").Value = 13138.652
$ws.Range("I74").Value = 13509
$ws.Range("K74").Value = 13509
$ws.Range("M74").Value = -12573
$ws.Range("H77").Value = 13138.652
$ws.Range("I77").Value = 13509
$ws.Range("K77").Value = 67545
$ws.Range("M77").Value = -62865
$ws.Range("H137").Value = 4168951.5
$ws.Range("I137").Value = 8335250
$ws.Range("J137").Value = 2653.1667
$ws.Range("K137").Value = 25005750
$ws.Range("L137").Value = 7959.500100000001
$ws.Range("M137").Value = -25003200
$ws.Range("N137").Value = -13059.5001
$ws.Range("H138").Value = 5297.683
$ws.Range("J138").Value = 6132.2812
$ws.Range("L138").Value = 18396.8436
$ws.Range("N138").Value = -28676.8436

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2403.1875
$ws.Range("I2").Value = 2336.5
$ws.Range("J2").Value = 2514.3333
$ws.Range("K2").Value = 2336.5
$ws.Range("L2").Value = 2514.3333
$ws.Range("M2").Value = -2223.5
$ws.Range("N2").Value = -2740.3333
$ws.Range("H32").Value = 2290174.5
$ws.Range("I32").Value = 1031039.2
$ws.Range("K32").Value = 1031039.2
$ws.Range("M32").Value = -1030752.2
$ws.Range("H61").Value = 8446
$ws.Range("I61").Value = 8446
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 8446
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -8234
$ws.Range("N61").ClearContents()
$ws.Range("H102").Value = 2863
$ws.Range("I102").Value = 1999
$ws.Range("K102").Value = 1999
$ws.Range("M102").Value = -377
$ws.Range("H110").Value = 576.75
$ws.Range("I110").Value = 482.8
$ws.Range("J110").Value = 733.3333
$ws.Range("K110").Value = 482.8
$ws.Range("L110").Value = 733.3333
$ws.Range("M110").Value = 1562.2
$ws.Range("N110").Value = -4823.3333
$ws.Range("H116").Value = 2403.1875
$ws.Range("I116").Value = 2336.5
$ws.Range("J116").Value = 2514.3333
$ws.Range("K116").Value = 2336.5
$ws.Range("L116").Value = 2514.3333
$ws.Range("M116").Value = -42.5
$ws.Range("N116").Value = -7102.3333
$ws.Range("H136").Value = 8446
$ws.Range("I136").Value = 8446
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 25338
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -22788
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2403.1875
$ws.Range("I3").Value = 2336.5
$ws.Range("J3").Value = 2514.3333
$ws.Range("K3").Value = 2336.5
$ws.Range("L3").Value = 2514.3333
$ws.Range("M3").Value = -2222.5
$ws.Range("N3").Value = -2742.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2239390.5
$ws.Range("I31").Value = 2085.0908
$ws.Range("J31").Value = 2786287.2
$ws.Range("K31").Value = 2085.0908
$ws.Range("L31").Value = 2786287.2
$ws.Range("M31").Value = -1790.0908
$ws.Range("N31").Value = -2786877.2
$ws.Range("H34").Value = 2239390.5
$ws.Range("I34").Value = 2085.0908
$ws.Range("J34").Value = 2786287.2
$ws.Range("K34").Value = 2085.0908
$ws.Range("L34").Value = 2786287.2
$ws.Range("M34").Value = -1883.0908
$ws.Range("N34").Value = -2786691.2
$ws.Range("H76").Value = 4998
$ws.Range("I76").Value = 4998
$ws.Range("K76").Value = 4998
$ws.Range("M76").Value = -4683
$ws.Range("H79").Value = 4998
$ws.Range("I79").Value = 4998
$ws.Range("K79").Value = 4998
$ws.Range("M79").Value = -3906
$ws.Range("H105").Value = 1971.4546
$ws.Range("I105").Value = 2076.4443
$ws.Range("K105").Value = 2076.4443
$ws.Range("M105").Value = -329.4443000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 8339192.5
$ws.Range("I68").Value = 1155.6666
$ws.Range("J68").Value = 11118538
$ws.Range("K68").Value = 3466.9998
$ws.Range("L68").Value = 33355614
$ws.Range("M68").Value = -2655.9998
$ws.Range("N68").Value = -33357236
$ws.Range("H71").Value = 8339192.5
$ws.Range("I71").Value = 1155.6666
$ws.Range("J71").Value = 11118538
$ws.Range("K71").Value = 10400.9994
$ws.Range("L71").Value = 100066842
$ws.Range("M71").Value = -6344.999400000001
$ws.Range("N71").Value = -100074954
$ws.Range("H114").Value = 3336.8462
$ws.Range("I114").Value = 2978.25
$ws.Range("J114").Value = 3496.2222
$ws.Range("K114").Value = 8934.75
$ws.Range("L114").Value = 10488.6666
$ws.Range("M114").Value = -5680.75
$ws.Range("N114").Value = -16996.6666
$ws.Range("H129").Value = 210
$ws.Range("I129").Value = 210
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 630
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 4370
$ws.Range("N129").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 30000000
$ws.Range("J103").Value = 30000000
$ws.Range("L103").Value = 30000000
$ws.Range("N103").Value = -30002344
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H126").Value = 5104.0527
$ws.Range("J126").Value = 6785.909
$ws.Range("L126").Value = 20357.727
$ws.Range("N126").Value = -25297.727
$ws.Range("H135").Value = 69999.14
$ws.Range("J135").Value = 69999.14
$ws.Range("L135").Value = 69999.14
$ws.Range("N135").Value = -80139.14

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1965.6364
$ws.Range("I16").Value = 625
$ws.Range("J16").Value = 7998.5
$ws.Range("K16").Value = 625
$ws.Range("L16").Value = 7998.5
$ws.Range("M16").Value = -455
$ws.Range("N16").Value = -8338.5
$ws.Range("H46").Value = 1624.25
$ws.Range("I46").Value = 833.6667
$ws.Range("K46").Value = 833.6667
$ws.Range("M46").Value = -645.6667
$ws.Range("H61").Value = 978.6875
$ws.Range("I61").Value = 989.1539
$ws.Range("J61").Value = 933.3333
$ws.Range("K61").Value = 989.1539
$ws.Range("L61").Value = 933.3333
$ws.Range("M61").Value = -787.1539
$ws.Range("N61").Value = -1337.3333
$ws.Range("H113").Value = 978.6875
$ws.Range("I113").Value = 989.1539
$ws.Range("J113").Value = 933.3333
$ws.Range("K113").Value = 989.1539
$ws.Range("L113").Value = 933.3333
$ws.Range("M113").Value = 1180.8461
$ws.Range("N113").Value = -5273.3333
$ws.Range("H138").Value = 94999
$ws.Range("J138").Value = 94999
$ws.Range("L138").Value = 94999
$ws.Range("N138").Value = -105279

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 20000
$ws.Range("I21").Value = 20000
$ws.Range("K21").Value = 20000
$ws.Range("M21").Value = -19765
$ws.Range("H35").Value = 20000
$ws.Range("I35").Value = 20000
$ws.Range("K35").Value = 20000
$ws.Range("M35").Value = -19710
$ws.Range("H40").Value = 3336008.2
$ws.Range("I40").Value = 2504012.5
$ws.Range("J40").Value = 5000000
$ws.Range("K40").Value = 2504012.5
$ws.Range("L40").Value = 5000000
$ws.Range("M40").Value = -2503863.5
$ws.Range("N40").Value = -5000298
$ws.Range("H113").Value = 433.64285
$ws.Range("I113").Value = 390.18182
$ws.Range("K113").Value = 1170.54546
$ws.Range("M113").Value = 999.45454
$ws.Range("H126").Value = 8163.375
$ws.Range("I126").Value = 10419.917
$ws.Range("K126").Value = 31259.751
$ws.Range("M126").Value = -28789.751
